$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + 1h volume change), and a rank swap
# between TrustWalletToken (row 46) and FTXToken (row 47), per commit diff.
# Price cells that look like plain decimal numbers (single "." separator)
# are forced to text format first so Excel keeps them as literal strings
# (matching the source data, e.g. preserving "21.00" instead of becoming 21).

$ws.Range("D2").Value = "37.746.81"
$ws.Range("E2").Value = "  -1.20%  "
$ws.Range("D3").Value = "2.032.34"
$ws.Range("E3").Value = "  -1.30%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.52"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.608"
$ws.Range("E6").Value = "  -1.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.09"
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.377"
$ws.Range("E9").Value = "  -2.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0822"
$ws.Range("E10").Value = "  +1.10%  "
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("D12").Value = "2.334.16"
$ws.Range("E12").Value = "  -1.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.43"
$ws.Range("E13").Value = "  -2.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.00"
$ws.Range("E14").Value = "  -1.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.758"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.17"
$ws.Range("E16").Value = "  -2.39%  "
$ws.Range("D17").Value = "2.025.99"
$ws.Range("E17").Value = "  -2.51%  "
$ws.Range("D18").Value = "37.722.16"
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.59"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.93"
$ws.Range("E20").Value = "  -5.75%  "
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.49"
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.42"
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.21"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.25"
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("E28").Value = "  -2.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.79"
$ws.Range("E29").Value = "  -1.48%  "
$ws.Range("E30").Value = "  -4.20%  "
$ws.Range("E31").Value = "  +0.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.19"
$ws.Range("E32").Value = "  +6.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.39"
$ws.Range("E33").Value = "  -4.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0603"
$ws.Range("E34").Value = "  -0.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.47"
$ws.Range("E35").Value = "  -2.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.31"
$ws.Range("E36").Value = "  +2.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.27"
$ws.Range("E37").Value = "  -2.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.32"
$ws.Range("E38").Value = "  +0.77%  "
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.78"
$ws.Range("E40").Value = "  +3.86%  "
$ws.Range("D41").Value = "1.538.75"
$ws.Range("E41").Value = "  +0.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0216"
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "95.70"
$ws.Range("E43").Value = "  -2.51%  "
$ws.Range("E44").Value = "  -2.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0911"
$ws.Range("E45").Value = "  -1.31%  "
$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.19"
$ws.Range("E46").Value = "  +3.08%  "
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.11"
$ws.Range("E47").Value = "  -2.71%  "
$ws.Range("E48").Value = "  -1.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.14"
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("E50").Value = "  -0.65%  "
$ws.Range("D51").Value = "2.223.94"
$ws.Range("E51").Value = "  -1.23%  "
